$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value that was bumped by one
# day (45188 -> 45189) for every data row (rows 2 through 295).
$ws.Range("C2:C295").Value = 45189
